$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-15 from 45183 to 45184
$ws.Range("C2:C15").Value = 45184
